# Ridership run on 20161026.
# Update the "Riders" (column C) and "Average" (column D) figures on the
# "Ridership" worksheet for August 2016 with the refreshed ridership counts.
# The line chart on this sheet is bound to Ridership!$C$2:$C$32 and
# Ridership!$D$2:$D$32, so Excel will refresh the chart cache automatically
# when these cells are recalculated/saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 180
$ws.Range("D2").Value = 92.49

$ws.Range("C3").Value = 171
$ws.Range("D3").Value = 93.59

$ws.Range("C4").Value = 209
$ws.Range("D4").Value = 100.35

$ws.Range("C5").Value = 168
$ws.Range("D5").Value = 97.48

$ws.Range("C6").Value = 184
$ws.Range("D6").Value = 94.38

$ws.Range("C9").Value = 139
$ws.Range("D9").Value = 93.55

$ws.Range("C10").Value = 219
$ws.Range("D10").Value = 96.1

$ws.Range("C11").Value = 213
$ws.Range("D11").Value = 102.65

$ws.Range("C12").Value = 195
$ws.Range("D12").Value = 99.47

$ws.Range("C13").Value = 249
$ws.Range("D13").Value = 97.74

$ws.Range("C14").Value = 100
$ws.Range("D14").Value = 41.64

$ws.Range("C15").Value = 72
$ws.Range("D15").Value = 33.98

$ws.Range("C16").Value = 182
$ws.Range("D16").Value = 95.51

$ws.Range("C17").Value = 184

$ws.Range("C18").Value = 177
$ws.Range("D18").Value = 104.14

$ws.Range("C19").Value = 236
$ws.Range("D19").Value = 102.2

$ws.Range("C20").Value = 201
$ws.Range("D20").Value = 99.94

$ws.Range("C21").Value = 87
$ws.Range("D21").Value = 42.58

$ws.Range("D22").Value = 35

$ws.Range("C23").Value = 185
$ws.Range("D23").Value = 97.46

$ws.Range("C24").Value = 198
$ws.Range("D24").Value = 99.75

$ws.Range("C25").Value = 191
$ws.Range("D25").Value = 105.84

$ws.Range("C26").Value = 206
$ws.Range("D26").Value = 104.24

$ws.Range("C27").Value = 189
$ws.Range("D27").Value = 101.79

$ws.Range("C28").Value = 303
$ws.Range("D28").Value = 47.9

$ws.Range("C29").Value = 118
$ws.Range("D29").Value = 36.66

$ws.Range("C30").Value = 229
$ws.Range("D30").Value = 100.26

$ws.Range("C31").Value = 251
$ws.Range("D31").Value = 102.6

$ws.Range("C32").Value = 173
$ws.Range("D32").Value = 107.13

$wb.Save()
